$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.922.73"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.641.96"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'215.38"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'0.5089"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.06395"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "'0.07771"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "'4.307"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "1.652.16"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "'0.5458"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "0.0₅7854"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "'64.66"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "25.984.85"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "'197.81"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").Value = "'4.443"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").Value = "'9.976"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").Value = "'6.035"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").Value = "'1.008"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").Value = "'1.879"
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("D25").Value = "'140.80"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").Value = "'6.909"
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("D28").Value = "'15.72"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").Value = "'1.241"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'0.05016"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").Value = "'3.263"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "'1.542"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").Value = "'2.365"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "'0.8944"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").Value = "'2.590"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("D37").Value = "1.135.83"
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("D38").Value = "'0.5515"
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("D39").Value = "'0.01557"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("B40").Value = "mCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D40").Value = "'2.557"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.006"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").Value = "'5.635"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "0.0₈128"
$ws.Range("E43").Value = "  +10.32%  "
$ws.Range("D44").Value = "'0.8173"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D45").Value = "'99.80"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("D46").Value = "1.780.80"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").Value = "'0.4531"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "'1.003"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").Value = "'54.96"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").Value = "'0.05088"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.09555"
$ws.Range("E51").Value = "  +3.17%  "
